$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$zoom = "[Zoom](https://nhh.zoom.us/j/66295455298?pwd=bWNBd2NlUE9PLzRCNHZQUEVuaVE1UT09 )."

$ws.Range("D3").Value = "02.09: Oppgaveseminar Aud Max/$zoom. Se \@ref(seminar) for oppgaver."
$ws.Range("D5").Value = "16.09: Oppgaveseminar Aud Max/$zoom. Se \@ref(seminar) for oppgaver."
$ws.Range("D7").Value = "30.09: Oppgaveseminar Aud Max/$zoom. Se \@ref(seminar) for oppgaver."
$ws.Range("D9").Value = "14.10: **Oversiktsforelesning: Hypotesetesting** i Aud Max/$zoom."
$ws.Range("D10").Value = "21.10: Oppgaveseminar Aud Max/$zoom. Se \@ref(seminar) for oppgaver."
$ws.Range("D11").Value = "28.10: **Oversiktsforelesning: Regresjon** i Aud Max/$zoom."
$ws.Range("D12").Value = "4.11: **Oversiktsforelesning: Tidsrekker** i Aud Max/$zoom."
$ws.Range("D13").Value = "11.11 Oppgaveseminar Aud Max/$zoom. Se \@ref(seminar) for oppgaver."

$ws.Range("D14").Select()
